# Fixed bug in path recognition
# Rebuilds the weekly schedule sheet: shifts the header class-column order,
# replaces the Jan date range with a longer Feb date range, and re-files
# each lesson entry into the (now single) "C" lesson column per day block
# with its own highlight colour, instead of the old scattered columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start clean: drop all merges + all cell content/formatting ---
$ws.Cells.UnMerge()
$ws.Cells.Clear()

$nl = [char]10

# --- header row (class labels), now reordered 6E,7E,8E,9E,10E,11E ---
$headers = @("6E","7E","8E","9E","10E","11E")
$headerCols = @("C","D","E","F","G","H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# --- week blocks: each one is a 6-row group (day header + 5 lesson-number rows)
#     followed by a 1-row grey separator, except the very last block. ---
$weeks = @(
    @{ Row = 2;  Label = "Среда${nl}15 Фев.";        Lessons = @{ 4 = @{ Col = "C"; Text = "История (Настя)"; Color = 0xee7bbb; Height = 25 };
                                                                    6 = @{ Col = "C"; Text = "Физ-ра (ЖК)";      Color = 0x0082c4; Height = 18 } } },
    @{ Row = 8;  Label = "Четверг${nl}16 Фев.";       Lessons = @{ 9  = @{ Col = "C"; Text = "География (ГН)";   Color = 0x00826a; Height = 23 };
                                                                    10 = @{ Col = "C"; Text = "Физика (ЮН)";      Color = 0x0085af; Height = 18 };
                                                                    12 = @{ Col = "C"; Text = "Информатика (Ок)"; Color = 0x0083dc; Height = 27 } } },
    @{ Row = 14; Label = "Пяница${nl}17 Фев.";        Lessons = @{} },
    @{ Row = 20; Label = "Суббота${nl}18 Фев.";       Lessons = @{} },
    @{ Row = 26; Label = "Понедельн.${nl}19 Фев.";    Lessons = @{} },
    @{ Row = 32; Label = "Понедельн.${nl}20 Фев.";    Lessons = @{} },
    @{ Row = 38; Label = "Вторник${nl}21 Фев.";       Lessons = @{} },
    @{ Row = 44; Label = "Среда${nl}22 Фев.";         Lessons = @{} },
    @{ Row = 50; Label = "Четверг${nl}23 Фев.";       Lessons = @{} },
    @{ Row = 56; Label = "Пяница${nl}24 Фев.";        Lessons = @{} },
    @{ Row = 62; Label = "Суббота${nl}25 Фев.";       Lessons = @{} },
    @{ Row = 68; Label = "Понедельн.${nl}26 Фев.";    Lessons = @{} },
    @{ Row = 74; Label = "Понедельн.${nl}27 Фев.";    Lessons = @{} },
    @{ Row = 80; Label = "Вторник${nl}28 Фев.";       Lessons = @{} }
)

foreach ($week in $weeks) {
    $r0 = $week.Row

    # day-label cell, merged A(r0):A(r0+4)
    $dayCell = $ws.Range("A" + $r0)
    $dayCell.Value = $week.Label
    $dayCell.HorizontalAlignment = -4108
    $dayCell.VerticalAlignment = -4108
    $dayCell.WrapText = $true
    $ws.Range("A" + $r0 + ":A" + ($r0 + 4)).Merge() | Out-Null

    for ($i = 0; $i -lt 5; $i++) {
        $r = $r0 + $i
        $numCell = $ws.Range("B" + $r)
        $numCell.Value = $i + 1
        $numCell.HorizontalAlignment = -4108
        $numCell.WrapText = $true

        if ($week.Lessons.ContainsKey($r)) {
            $lesson = $week.Lessons[$r]
            $lessonCell = $ws.Range($lesson.Col + $r)
            $lessonCell.Value = $lesson.Text
            $lessonCell.VerticalAlignment = -4108
            $lessonCell.WrapText = $true
            $lessonCell.Interior.Color = $lesson.Color
            $ws.Rows.Item($r).RowHeight = $lesson.Height
        }
    }

    # grey separator row right after the block, merged across A:H - skip after the final week
    $sepRow = $r0 + 5
    if ($sepRow -le 79) {
        $ws.Range("A" + $sepRow + ":H" + $sepRow).Merge() | Out-Null
        $ws.Range("A" + $sepRow).Interior.Color = 0x009E9E9E
    }
}
